# The sheet's data rows (93..200) shift down by one to make room for a new
# weekly price record. A new row 93 is inserted, and all rows that were at
# 93..200 move to 94..201 (this is a native Excel row insert, so formats -
# including the date style on column D - carry over automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("93").Insert()

# Populate the newly inserted row 93 with the new weekly record.
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 44638
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = 100112043
$ws.Range("G93").Value = "Pepino dulce"
$ws.Range("H93").Value = "Cultivar IV Región"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 140
$ws.Range("K93").Value = 16000
$ws.Range("L93").Value = 16000
$ws.Range("M93").Value = 16000
$ws.Range("N93").Value = "$/bandeja 18 kilos"
$ws.Range("O93").Value = "Provincia de Limarí"
$ws.Range("P93").Value = 889
$ws.Range("Q93").Value = 18
$ws.Range("R93").Value = "Hortaliza"
